$wb = $excel.ActiveWorkbook

# --- Sheet "Precios" (sheet2.xml): new DB design adds Building.Code / Building.Name columns ---
$ws = $wb.Worksheets.Item("Precios")

# Insert two new columns at B:C, shifting the old B..N block to D..P
$ws.Columns("B:C").EntireColumn.Insert()

# New header row (row 1): merged "Edificio" header over B1:C1 (style matches neighboring s=3 header cells)
$ws.Range("B1").Value = "Edificio"

# New sub-header row (row 2): Building.Code / Building.Name
$ws.Range("B2").Value = "Building.Code"
$ws.Range("C2").Value = "Building.Name"

# Column widths for the two new columns
$ws.Columns("B").ColumnWidth = 16.5
$ws.Columns("C").ColumnWidth = 16.5

# Refresh the AutoFilter over the new, wider range
$ws.AutoFilterMode = $false
$ws.Range("A2:M3").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name for Precios in sync with the new range
$wb.Names.Item("Precios!_FilterDatabase").RefersTo = "=Precios!`$A`$2:`$M`$3"

# --- Active tab moves from "Recursos" to "Precios" ---
$ws.Activate()
